$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column/row D for item "D" in the AHP pairwise comparison matrix
$ws.Range("D1").Value = "D"

$ws.Range("A4").Value = "D"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 1

$ws.Range("D2").Value = 0.2
$ws.Range("D3").Value = 0.2

# Remove the explicit number format style that was on C2 (now plain General)
$ws.Range("C2").ClearFormats()

$ws.Range("D4").Select()
